# Regenerate the whole experiment data on sheet "Arkusz1" (A1:E35).
# Column A = graph_id, B = size, C = no_comms, D = inside_prob, E = outside_prob.
# Row 1 (headers) is left untouched; rows 2-35 get fresh generated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# (graph_id, size, no_comms, inside_prob, outside_prob) for rows 2-35.
# The last row's outside_prob (E35) is left alone on purpose: it was
# already stored as the text "0.08" before this edit and stays exactly
# that way afterwards too, so it is skipped below.
$data = @(
    @(0, 50, 2, 0.4, 0.01),
    @(1, 50, 2, 0.3, 0.02),
    @(2, 50, 2, 0.5, 0.05),
    @(3, 50, 2, 0.4, 0.08),
    @(4, 50, 2, 0.4, 0.1),
    @(5, 50, 2, 0.7, 0.12),
    @(6, 50, 3, 0.35, 0.06),
    @(7, 50, 3, 0.3, 0.02),
    @(8, 50, 3, 0.4, 0.05),
    @(9, 50, 3, 0.4, 0.02),
    @(10, 50, 3, 0.4, 0.05),
    @(11, 50, 3, 0.4, 0.08),
    @(12, 50, 4, 0.4, 0.01),
    @(13, 50, 4, 0.3, 0.03),
    @(14, 50, 4, 0.5, 0.02),
    @(15, 50, 4, 0.3, 0.02),
    @(16, 50, 4, 0.3, 0.02),
    @(17, 50, 4, 0.35, 0.02),
    @(18, 50, 4, 0.5, 0.07),
    @(19, 50, 5, 0.5, 0.05),
    @(20, 50, 5, 0.4, 0.01),
    @(21, 50, 5, 0.3, 0.02),
    @(22, 50, 5, 0.3, 0.02),
    @(23, 50, 5, 0.4, 0.05),
    @(24, 50, 5, 0.5, 0.08),
    @(25, 50, 6, 0.45, 0.1),
    @(26, 50, 6, 0.4, 0.01),
    @(27, 50, 6, 0.5, 0.02),
    @(28, 50, 6, 0.4, 0.05),
    @(29, 50, 6, 0.4, 0.03),
    @(30, 50, 6, 0.8, 0.06),
    @(31, 50, 6, 0.65, 0.08),
    @(32, 50, 6, 0.6, 0.05),
    @(33, 50, 6, 0.55000000000000004, $null)
)

$row = 2
foreach ($record in $data) {
    $ws.Cells.Item($row, 1).Value = $record[0]
    $ws.Cells.Item($row, 2).Value = $record[1]
    $ws.Cells.Item($row, 3).Value = $record[2]
    $ws.Cells.Item($row, 4).Value = $record[3]
    if ($null -ne $record[4]) {
        $ws.Cells.Item($row, 5).Value = $record[4]
    }
    $row = $row + 1
}
